# Update batting-performance counts (ODI and T20 appearance totals)
# for several players on the active sheet, reflecting additional
# scraped match data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "D2"  = 24
    "E2"  = 36
    "D8"  = 66
    "E8"  = 38
    "D9"  = 39
    "E9"  = 57
    "D15" = 4
    "D22" = 23
    "E22" = 42
    "D24" = 19
    "E27" = 7
    "E28" = 11
    "D30" = 220
    "D31" = 95
    "E31" = 54
    "D36" = 23
    "E36" = 41
    "E39" = 62
    "D43" = 18
    "E43" = 17
    "D46" = 50
    "E46" = 87
    "E47" = 37
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
